# Reporting year for this workbook changed from 2019 to 2015:
# the sheet now carries one data row recording the new year value,
# and the selection moves to the next empty row (A3) as it would
# after Excel's "Sort/Filter Database" range grows by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row: year = 2015 in column A of row 2.
$ws.Range("A2").Value = 2015

# Move the active selection to A3, the next empty row below the data.
[void]$ws.Range("A3").Select()
